$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Karsten")

# New log entries added on the "Karsten" sheet (rows 6 and 7): times worked
# plus a description of the work done, pulled from the shared-string pool.
$ws1.Range("A6").Value = 0.57638888888888895
$ws1.Range("B6").Value = 0.625
$ws1.Range("D6").Value = "Douwe's code bekeken en in apart vsC bestand gekeken of het beter kon"

$ws1.Range("A7").Value = 0.4375
$ws1.Range("B7").Value = 0.54166666666666663
$ws1.Range("D7").Value = "Verschillende defs en commands aangepast voor performance en Ease of Use"

# The "Karsten" tab becomes the active/selected sheet (previously "Douwe"
# was active), with D9 as the selected cell.
[void]$ws1.Select()
[void]$ws1.Range("D9").Select()
